# EIA Table A.5.A - Relative Standard Error for Net Generation by Fuel Type
# Monthly refresh: October 2016 data -> November 2016 data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the report subtitle (row 2) to reflect the new reporting month
$ws.Range("A2").Value = "Industrial Sector by Census Division and State, November 2016"

# Refresh the relative standard error figures for each Census division / state row
# Row 4
$ws.Range("B4").Value = 81
$ws.Range("C4").Value = 141
$ws.Range("E4").Value = 36
$ws.Range("H4").Value = 47
# Row 5
$ws.Range("C5").Value = 436
$ws.Range("E5").Value = 55
# Row 6
$ws.Range("C6").Value = 128
$ws.Range("E6").Value = 51
$ws.Range("H6").Value = 47
# Row 7
$ws.Range("B7").Value = 195
$ws.Range("C7").Value = 4303
$ws.Range("H7").Value = 788
# Row 8
$ws.Range("C8").Value = 418
$ws.Range("E8").Value = 209
# Row 9
$ws.Range("B9").Value = 22
$ws.Range("D9").Value = 69
$ws.Range("E9").Value = 24
$ws.Range("F9").Value = 29
$ws.Range("H9").Value = 182
# Row 10
$ws.Range("C10").Value = 529
$ws.Range("D10").Value = 128
$ws.Range("E10").Value = 64
$ws.Range("F10").Value = 75
# Row 11
$ws.Range("C11").Value = 46
$ws.Range("E11").Value = 38
$ws.Range("H11").Value = 182
# Row 12
$ws.Range("B12").Value = 47
$ws.Range("C12").Value = 25
$ws.Range("D12").Value = 81
$ws.Range("E12").Value = 31
$ws.Range("F12").Value = 28
# Row 13
$ws.Range("B13").Value = 9
$ws.Range("C13").Value = 23
$ws.Range("D13").Value = 78
$ws.Range("E13").Value = 18
$ws.Range("F13").Value = 21
$ws.Range("H13").Value = 82
# Row 14
$ws.Range("B14").Value = 8
$ws.Range("E14").Value = 50
$ws.Range("F14").Value = 99
# Row 15
$ws.Range("B15").Value = 614
$ws.Range("C15").Value = 8
$ws.Range("E15").Value = 32
$ws.Range("F15").Value = 16
# Row 16
$ws.Range("B16").Value = 99
$ws.Range("C16").Value = 19
$ws.Range("D16").Value = 89
$ws.Range("E16").Value = 33
$ws.Range("H16").Value = 208
# Row 17
$ws.Range("B17").Value = 125
$ws.Range("C17").Value = 213
$ws.Range("D17").Value = 339
$ws.Range("E17").Value = 43
$ws.Range("F17").Value = 168
# Row 18
$ws.Range("B18").Value = 18
$ws.Range("C18").Value = 374
$ws.Range("E18").Value = 37
$ws.Range("H18").Value = 89
# Row 19
$ws.Range("B19").Value = 15
$ws.Range("C19").Value = 224
$ws.Range("D19").Value = 168
$ws.Range("E19").Value = 29
$ws.Range("F19").Value = 100
$ws.Range("H19").Value = 97
# Row 20
$ws.Range("B20").Value = 15
$ws.Range("C20").Value = 429
$ws.Range("D20").Value = 168
$ws.Range("E20").Value = 30
# Row 21
$ws.Range("E21").Value = 112
# Row 22
$ws.Range("B22").Value = 37
$ws.Range("C22").Value = 409
$ws.Range("E22").Value = 87
$ws.Range("H22").Value = 97
# Row 23
$ws.Range("B23").Value = 175
$ws.Range("E23").Value = 306
# Row 24
$ws.Range("B24").Value = 44
$ws.Range("E24").Value = 457
# Row 25
$ws.Range("B25").Value = 90
$ws.Range("C25").Value = 290
$ws.Range("E25").Value = 175
$ws.Range("F25").Value = 100
# Row 26
$ws.Range("B26").Value = 21
$ws.Range("C26").Value = 74
$ws.Range("E26").Value = 9
$ws.Range("H26").Value = 46
# Row 28
$ws.Range("B28").Value = 85
$ws.Range("C28").Value = 201
# Row 29
$ws.Range("B29").Value = 42
$ws.Range("C29").Value = 78
$ws.Range("E29").Value = 25
$ws.Range("H29").Value = 288
# Row 30
$ws.Range("C30").Value = 277
$ws.Range("E30").Value = 87
# Row 31
$ws.Range("B31").Value = 82
$ws.Range("C31").Value = 666
$ws.Range("E31").Value = 70
$ws.Range("H31").Value = 1086
# Row 32
$ws.Range("B32").Value = 9
$ws.Range("C32").Value = 8
$ws.Range("E32").Value = 75
# Row 33
$ws.Range("B33").Value = 25
$ws.Range("C33").Value = 570
$ws.Range("E33").Value = 29
$ws.Range("H33").Value = 405
# Row 34
$ws.Range("H34").Value = 30
# Row 35
$ws.Range("B35").Value = 6
$ws.Range("C35").Value = 88
$ws.Range("E35").Value = 15
$ws.Range("F35").Value = 72
# Row 36
$ws.Range("B36").Value = 65
$ws.Range("C36").Value = 104
$ws.Range("E36").Value = 22
$ws.Range("F36").Value = 119
# Row 37
$ws.Range("E37").Value = 71
# Row 38
$ws.Range("E38").Value = 38
# Row 39
$ws.Range("C39").Value = 148
# Row 40
$ws.Range("B40").Value = 42
$ws.Range("C40").Value = 67
$ws.Range("D40").Value = 55
$ws.Range("F40").Value = 7
# Row 41
$ws.Range("C41").Value = 49
$ws.Range("E41").Value = 24
# Row 42
$ws.Range("D42").Value = 78
$ws.Range("F42").Value = 7
# Row 43
$ws.Range("B43").Value = 50
$ws.Range("C43").Value = 88
$ws.Range("E43").Value = 84
# Row 44
$ws.Range("C44").Value = 408
$ws.Range("D44").Value = 59
$ws.Range("F44").Value = 12
# Row 45
$ws.Range("B45").Value = 42
$ws.Range("C45").Value = 743
$ws.Range("F45").Value = 9
# Row 46
$ws.Range("B46").Value = 370
$ws.Range("C46").Value = 678
$ws.Range("E46").Value = 80
# Row 47
$ws.Range("B47").Value = 105
$ws.Range("E47").Value = 50
# Row 48
$ws.Range("B48").Value = 308
# Row 49
$ws.Range("E49").Value = 21
# Row 50
$ws.Range("C50").Value = 3020
$ws.Range("E50").Value = 0
# Row 51
$ws.Range("C51").Value = 1341
$ws.Range("E51").Value = 16
$ws.Range("F51").Value = 457
# Row 52
$ws.Range("B52").Value = 46
$ws.Range("C52").Value = 372
$ws.Range("E52").Value = 16
$ws.Range("F52").Value = 7
# Row 53
$ws.Range("C53").Value = 95
$ws.Range("E53").Value = 3
$ws.Range("F53").Value = 8
# Row 54
$ws.Range("C54").Value = 405
$ws.Range("E54").Value = 3
$ws.Range("F54").Value = 8
# Row 55
$ws.Range("E55").Value = 72
# Row 56
$ws.Range("C56").Value = 82
# Row 57
$ws.Range("B57").Value = 201
$ws.Range("C57").Value = 53
$ws.Range("E57").Value = 106
$ws.Range("F57").Value = 130
$ws.Range("H57").Value = 149
# Row 58
$ws.Range("C58").Value = 19
$ws.Range("E58").Value = 106
# Row 59
$ws.Range("B59").Value = 201
$ws.Range("C59").Value = 62
$ws.Range("F59").Value = 130
$ws.Range("H59").Value = 149
# Row 60
$ws.Range("B60").Value = 7
$ws.Range("C60").Value = 33
$ws.Range("D60").Value = 34
$ws.Range("F60").Value = 7
$ws.Range("H60").Value = 32
